$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

function Set-TranslationRow($row, $textId, $typography, $alignment, $gb, $direction) {
    $ws.Range("B$row").Value = $textId
    $ws.Range("C$row").Value = $typography
    $ws.Range("D$row").Value = $alignment
    $ws.Range("E$row").Value = $gb
    $ws.Range("F$row").Value = $direction

    # Creating a value in a previously-empty row bakes the column's default
    # style index onto the new cell (s="1"); the target cells carry no
    # explicit style, so normalize the style back after writing the value.
    $ws.Range("B$row`:F$row").Style = "Normal"
}

# Row 36: new text field "White"
Set-TranslationRow 36 "SingleUseId42" "Medium" "Center" "White" "LTR"

# Row 37: new text field "Dark"
Set-TranslationRow 37 "SingleUseId43" "Medium" "Center" "Dark" "LTR"

# Row 38: new text field "Urban"
Set-TranslationRow 38 "SingleUseId44" "Medium" "Center" "Urban" "LTR"
